# Add a "Save" column (column H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: same look/feel (bold, bordered, centered) as the other header cells.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells for the new column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
